# Replace the inline placeholder picture (Submission Flowchart image) with a
# hyperlink whose visible text is the image's source URL.
$d = $word.ActiveDocument

$imageUrl = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/WCP.jpg"

$shape = $d.InlineShapes(1)
$shapeRange = $shape.Range

# Remove the picture itself (keeps the surrounding paragraph intact).
$shape.Delete()

# Insert a hyperlink in its place, displaying the URL as text.
$d.Hyperlinks.Add($shapeRange, $imageUrl, [Type]::Missing, [Type]::Missing, $imageUrl)
